# Update stats for 2026-02 (row 27)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B27").Value = 6546
$ws.Range("C27").Value = 1018
$ws.Range("E27").Value = 932.6567369385884
$ws.Range("F27").Value = 10.01680672268908
$ws.Range("G27").Value = 7.4973600844773
